$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold text values (coin price & 1h volume %, e.g. "59.407.95").
# Some prices are plain decimals ("562.75") that Excel would otherwise auto-
# detect as numbers; force those cells to Text first so the stored value stays
# a string (matching the rest of the sheet), then drop back to the Normal style
# so no stray per-cell formatting is left behind.

$ws.Range("D2").Value = "59.407.95"
$ws.Range("E2").Value = "  +3.43%  "
$ws.Range("D3").Value = "3.000.24"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.56%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.523"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "2.986.90"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "3.487.08"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.11%  "
$ws.Range("D18").Value = "2.990.26"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").Value = "59.212.95"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.720"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0990"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.992"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.72%  "
$ws.Range("E36").Value = "  +13.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "400.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.08%  "
$ws.Range("D42").Value = "2.759.12"
$ws.Range("E42").Value = "  +4.88%  "
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  +5.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +23.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.82%  "
